$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-40 (A column index 0-38) with new B/C values
$ws.Cells.Item(2, 2).Value = 300.5
$ws.Cells.Item(2, 3).Value = 704.3133693271919
$ws.Cells.Item(3, 2).Value = 314.5
$ws.Cells.Item(3, 3).Value = 686.9740557409187
$ws.Cells.Item(4, 2).Value = 328.4
$ws.Cells.Item(4, 3).Value = 697.1478495955141
$ws.Cells.Item(5, 2).Value = 342
$ws.Cells.Item(5, 3).Value = 708.8310419864907
$ws.Cells.Item(6, 2).Value = 356
$ws.Cells.Item(6, 3).Value = 712.5912987658855
$ws.Cells.Item(7, 2).Value = 369.6
$ws.Cells.Item(7, 3).Value = 718.2218431714959
$ws.Cells.Item(8, 2).Value = 383.5
$ws.Cells.Item(8, 3).Value = 726.018359510916
$ws.Cells.Item(9, 2).Value = 397.5
$ws.Cells.Item(9, 3).Value = 724.8644724312525
$ws.Cells.Item(10, 2).Value = 411.1
$ws.Cells.Item(10, 3).Value = 722.577858361384
$ws.Cells.Item(11, 2).Value = 425
$ws.Cells.Item(11, 3).Value = 727.0246601275287
$ws.Cells.Item(12, 2).Value = 439
$ws.Cells.Item(12, 3).Value = 731.5756397161532
$ws.Cells.Item(13, 2).Value = 452.6
$ws.Cells.Item(13, 3).Value = 731.7144591617586
$ws.Cells.Item(14, 2).Value = 466.6
$ws.Cells.Item(14, 3).Value = 732.8499370130581
$ws.Cells.Item(15, 2).Value = 480.5
$ws.Cells.Item(15, 3).Value = 732.0109699673338
$ws.Cells.Item(16, 2).Value = 494.1
$ws.Cells.Item(16, 3).Value = 733.5484968087336
$ws.Cells.Item(17, 2).Value = 508.3
$ws.Cells.Item(17, 3).Value = 730.5707453100218
$ws.Cells.Item(18, 2).Value = 522
$ws.Cells.Item(18, 3).Value = 735.8716051000247
$ws.Cells.Item(19, 2).Value = 535.5999999999999
$ws.Cells.Item(19, 3).Value = 737.6614333026721
$ws.Cells.Item(20, 2).Value = 550
$ws.Cells.Item(20, 3).Value = 738.4388434966693
$ws.Cells.Item(21, 2).Value = 563.5
$ws.Cells.Item(21, 3).Value = 739.7035098560823
$ws.Cells.Item(22, 2).Value = 588.4000000000001
$ws.Cells.Item(22, 3).Value = 740.5995569664814
$ws.Cells.Item(23, 2).Value = 614
$ws.Cells.Item(23, 3).Value = 741.3664004880527
$ws.Cells.Item(24, 2).Value = 639
$ws.Cells.Item(24, 3).Value = 741.8630918425429
$ws.Cells.Item(25, 2).Value = 664.6
$ws.Cells.Item(25, 3).Value = 742.4857380700089
$ws.Cells.Item(26, 2).Value = 689.5
$ws.Cells.Item(26, 3).Value = 742.6864641844055
$ws.Cells.Item(27, 2).Value = 715
$ws.Cells.Item(27, 3).Value = 742.9212700500453
$ws.Cells.Item(28, 2).Value = 740
$ws.Cells.Item(28, 3).Value = 743.1444607878257
$ws.Cells.Item(29, 2).Value = 765
$ws.Cells.Item(29, 3).Value = 743.3359770652181
$ws.Cells.Item(30, 2).Value = 790.5
$ws.Cells.Item(30, 3).Value = 743.5288691962638
$ws.Cells.Item(31, 2).Value = 816
$ws.Cells.Item(31, 3).Value = 743.5043738865943
$ws.Cells.Item(32, 2).Value = 841.3000000000001
$ws.Cells.Item(32, 3).Value = 743.5812612943682
$ws.Cells.Item(33, 2).Value = 866
$ws.Cells.Item(33, 3).Value = 743.6769156722312
$ws.Cells.Item(34, 2).Value = 891.5999999999999
$ws.Cells.Item(34, 3).Value = 743.7060042859711
$ws.Cells.Item(35, 2).Value = 917
$ws.Cells.Item(35, 3).Value = 743.7770749892936
$ws.Cells.Item(36, 2).Value = 942.4
$ws.Cells.Item(36, 3).Value = 743.8176213373504
$ws.Cells.Item(37, 2).Value = 967.3000000000001
$ws.Cells.Item(37, 3).Value = 743.7738191475862
$ws.Cells.Item(38, 2).Value = 992.7
$ws.Cells.Item(38, 3).Value = 743.7862621196833
$ws.Cells.Item(39, 2).Value = 1018
$ws.Cells.Item(39, 3).Value = 743.8876858174571
$ws.Cells.Item(40, 2).Value = 1055
$ws.Cells.Item(40, 3).Value = 743.9433101416038

# Add new rows 41-65 (A column index 39-63), copying style of A40 for column A
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(41, 1).PasteSpecial(-4122)
$ws.Cells.Item(41, 1).Value = 39
$ws.Cells.Item(41, 2).Value = 1083
$ws.Cells.Item(41, 3).Value = 744.1013491696901
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(42, 1).PasteSpecial(-4122)
$ws.Cells.Item(42, 1).Value = 40
$ws.Cells.Item(42, 2).Value = 1112
$ws.Cells.Item(42, 3).Value = 744.132675177922
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(43, 1).PasteSpecial(-4122)
$ws.Cells.Item(43, 1).Value = 41
$ws.Cells.Item(43, 2).Value = 1143
$ws.Cells.Item(43, 3).Value = 744.243805770095
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(44, 1).PasteSpecial(-4122)
$ws.Cells.Item(44, 1).Value = 42
$ws.Cells.Item(44, 2).Value = 1174
$ws.Cells.Item(44, 3).Value = 744.3707989247539
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(45, 1).PasteSpecial(-4122)
$ws.Cells.Item(45, 1).Value = 43
$ws.Cells.Item(45, 2).Value = 1205
$ws.Cells.Item(45, 3).Value = 744.56069659945
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(46, 1).PasteSpecial(-4122)
$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 2).Value = 1237
$ws.Cells.Item(46, 3).Value = 744.7778539927397
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(47, 1).PasteSpecial(-4122)
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 2).Value = 1271
$ws.Cells.Item(47, 3).Value = 745.0123155428034
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(48, 1).PasteSpecial(-4122)
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 2).Value = 1306
$ws.Cells.Item(48, 3).Value = 745.3141678235297
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(49, 1).PasteSpecial(-4122)
$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 2).Value = 1341
$ws.Cells.Item(49, 3).Value = 745.7356408565596
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(50, 1).PasteSpecial(-4122)
$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(50, 2).Value = 1377
$ws.Cells.Item(50, 3).Value = 746.2094573807976
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(51, 1).PasteSpecial(-4122)
$ws.Cells.Item(51, 1).Value = 49
$ws.Cells.Item(51, 2).Value = 1414
$ws.Cells.Item(51, 3).Value = 746.7674708976357
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(52, 1).PasteSpecial(-4122)
$ws.Cells.Item(52, 1).Value = 50
$ws.Cells.Item(52, 2).Value = 1452
$ws.Cells.Item(52, 3).Value = 747.4083864059677
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(53, 1).PasteSpecial(-4122)
$ws.Cells.Item(53, 1).Value = 51
$ws.Cells.Item(53, 2).Value = 1492
$ws.Cells.Item(53, 3).Value = 747.927753055433
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(54, 1).PasteSpecial(-4122)
$ws.Cells.Item(54, 1).Value = 52
$ws.Cells.Item(54, 2).Value = 1532
$ws.Cells.Item(54, 3).Value = 748.4336959086013
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(55, 1).PasteSpecial(-4122)
$ws.Cells.Item(55, 1).Value = 53
$ws.Cells.Item(55, 2).Value = 1573
$ws.Cells.Item(55, 3).Value = 748.6162162658738
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(56, 1).PasteSpecial(-4122)
$ws.Cells.Item(56, 1).Value = 54
$ws.Cells.Item(56, 2).Value = 1616
$ws.Cells.Item(56, 3).Value = 748.4099882123508
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(57, 1).PasteSpecial(-4122)
$ws.Cells.Item(57, 1).Value = 55
$ws.Cells.Item(57, 2).Value = 1659
$ws.Cells.Item(57, 3).Value = 748.1248621132863
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(58, 1).PasteSpecial(-4122)
$ws.Cells.Item(58, 1).Value = 56
$ws.Cells.Item(58, 2).Value = 1704
$ws.Cells.Item(58, 3).Value = 747.9217778639802
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(59, 1).PasteSpecial(-4122)
$ws.Cells.Item(59, 1).Value = 57
$ws.Cells.Item(59, 2).Value = 1751
$ws.Cells.Item(59, 3).Value = 747.7542444689344
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(60, 1).PasteSpecial(-4122)
$ws.Cells.Item(60, 1).Value = 58
$ws.Cells.Item(60, 2).Value = 2098
$ws.Cells.Item(60, 3).Value = 746.2247628140517
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(61, 1).PasteSpecial(-4122)
$ws.Cells.Item(61, 1).Value = 59
$ws.Cells.Item(61, 2).Value = 2938
$ws.Cells.Item(61, 3).Value = 743.9612005673002
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(62, 1).PasteSpecial(-4122)
$ws.Cells.Item(62, 1).Value = 60
$ws.Cells.Item(62, 2).Value = 3060
$ws.Cells.Item(62, 3).Value = 743.7711936857384
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(63, 1).PasteSpecial(-4122)
$ws.Cells.Item(63, 1).Value = 61
$ws.Cells.Item(63, 2).Value = 3184
$ws.Cells.Item(63, 3).Value = 743.3964595170428
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(64, 1).PasteSpecial(-4122)
$ws.Cells.Item(64, 1).Value = 62
$ws.Cells.Item(64, 2).Value = 3307
$ws.Cells.Item(64, 3).Value = 743.2959976072517
$ws.Cells.Item(40, 1).Copy()
$ws.Cells.Item(65, 1).PasteSpecial(-4122)
$ws.Cells.Item(65, 1).Value = 63
$ws.Cells.Item(65, 2).Value = 3430
$ws.Cells.Item(65, 3).Value = 742.8183544393452

$excel.CutCopyMode = 0
Write-Host "done"